$p = $ppt.ActivePresentation

# ----------------------------------------------------------------------
# Slides 13-15: underline the leading field-name token of various bullet
# paragraphs (splits the run so the field name gets u="sng").
# Each entry: slide index, shape index, paragraph index, # chars to underline
# (counted from the start of the paragraph).
# ----------------------------------------------------------------------

function Underline-FieldName {
    param($slideIndex, $shapeIndex, $paraIndex, $charCount)

    $slide = $p.Slides.Item($slideIndex)
    $shape = $slide.Shapes.Item($shapeIndex)
    $tr = $shape.TextFrame.TextRange
    $para = $tr.Paragraphs($paraIndex, 1)
    $chars = $tr.Characters($para.Start, $charCount)
    $chars.Font.Underline = $true
}

# --- Slide 13 : Company / Survey -------------------------------------
Underline-FieldName 13 2 2 2   # "ID"       (SMALLINT) - Auto Increment ID
Underline-FieldName 13 2 3 4   # "Name"     (VARCHAR) - Name of the company
Underline-FieldName 13 2 4 7   # "Address"  (VARCHAR) - Address of the company
Underline-FieldName 13 2 5 5   # "Phone"    (VARCHAR) - Phone number of contact
Underline-FieldName 13 2 7 2   # "ID"       (INT) - Auto Increment ID
Underline-FieldName 13 2 8 3   # "CID"      (SMALLINT) - Company ID ...
Underline-FieldName 13 2 9 9   # "StartDate"(DATE) - The start date for the survey
Underline-FieldName 13 2 10 7  # "EndDate"  (DATE) - The end date for the survey

# --- Slide 14 : Department / Participant ------------------------------
Underline-FieldName 14 2 2 2   # "ID"        (INT) - Auto Increment ID
Underline-FieldName 14 2 3 3   # "SID"       (INT) - Survey ID ...
Underline-FieldName 14 2 4 4   # "Name"      (VARCHAR) - Name of the company
Underline-FieldName 14 2 6 2   # "ID"        (INT) - Auto Increment ID
Underline-FieldName 14 2 7 3   # "DID"       (INT) - Department ID ...
Underline-FieldName 14 2 8 5   # "Email"     (VARCHAR) - The participant's email
Underline-FieldName 14 2 9 9   # "Submitted" (TINYINT) - Indicated submission status (...)

# --- Slide 15 : Question / SurveyQuestion / Response -------------------
Underline-FieldName 15 2 2 2   # "ID"       (SMALLINT) - Configurable Question ID
Underline-FieldName 15 2 3 2   # "LS"       (VARCHAR) - Left statement
Underline-FieldName 15 2 4 2   # "RS"       (VARCHAR) - Right statement
Underline-FieldName 15 2 6 3   # "SID"      (INT) - Survey ID ...
Underline-FieldName 15 2 7 3   # "QID"      (INT) - Question ID ...
Underline-FieldName 15 2 8 5   # "Order"    (SMALLINT) - The order of the question in the survey
Underline-FieldName 15 2 10 2  # "ID"       (INT) - Auto Increment ID
Underline-FieldName 15 2 11 3  # "DID"      (INT) - Department ID ...
Underline-FieldName 15 2 12 3  # "QID"      (SMALLINT) - Question ID ...
Underline-FieldName 15 2 13 8  # "Response" (TINYINT) - Participant response (1-6)

# ----------------------------------------------------------------------
# Slide 16: merge the two runs "Relationships maintain " + "cascading:"
# back into a single run "Relationships maintain cascading:" (same rPr).
# ----------------------------------------------------------------------

$slide16 = $p.Slides.Item(16)
$shape16 = $slide16.Shapes.Item(3)
$tr16 = $shape16.TextFrame.TextRange
$para1 = $tr16.Paragraphs(1, 1)
# Force an actual (unrelated) text change so the engine collapses the
# paragraph into a single run, then restore the final desired text.
$para1.Text = "zzz"
$tr16.Paragraphs(1, 1).Text = "Relationships maintain cascading:"
